# TORIBDA feat: added generateItemsDetail method
# Fill in the PDCA log entries for the "generateItemsDetails" row (row 8)
# and correct the "generateReceipt" row (row 7) Check/Action text, then
# move the selection the way the author left it when they saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 (generateReceipt) - Check / Action columns
$ws.Range("C7").Value = "7 mins 18 seconds"
$ws.Range("D7").Value = "have only created pseudocode and added the method for the sub-tasks all the while setting the receipt format"

# Row 8 (generateItemsDetails) - Do / Check / Action / Plan(notes) columns
$ws.Range("B8").Value = "3 mins"
$ws.Range("C8").Value = "12 mins 27 seconds"
$ws.Range("D8").Value = "Time spent more on formatting so that the PosMachineTest would result to success"
$ws.Range("E8").Value = "Check the required formatting rather than code blindly of whatever format (try not to overlook when there is a provided format)."

# Row 8 grew taller to fit the new wrapped text.
$ws.Rows.Item(8).RowHeight = 75

# Leave the selection where the author left it on save.
$ws.Range("D9").Select()
